$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "URL" column (D) with team slug/lookup values.
# Only the first 47 data rows received values in the source edit
# (the remaining rows were left unfilled, per the commit message
# "Fixed NBA binary for all sheets, some missing though").

$values = @{
    1  = "URL"
    2  = "None"
    3  = "None"
    4  = "akron"
    5  = "alabama"
    6  = "alabamaam"
    7  = "uab"
    8  = "alabamast"
    9  = "albanyst"
    10 = "alcornst"
    11 = "american"
    12 = "appalst"
    13 = "Arizona"
    14 = "arizonast"
    15 = "arkansas"
    16 = "arkansaslr"
    17 = "arkansaspb"
    18 = "arkansasst"
    19 = "army"
    20 = "auburn"
    21 = "austinpeay"
    22 = "ballst"
    23 = "baylor"
    24 = "belmont"
    25 = "bethcook"
    26 = "None"
    27 = "uab"
    28 = "boisest"
    29 = "bostoncoll"
    30 = "bostonuniv"
    31 = "bowlgreen"
    32 = "bradley"
    33 = "byu"
    34 = "brown"
    35 = "None"
    36 = "bucknell"
    37 = "buffalost"
    38 = "butler"
    39 = "calpoly"
    40 = "calstbake"
    41 = "calstfull"
    42 = "calstnorth"
    43 = "california"
    44 = "None"
    45 = "ucirvine"
    46 = "ucriver"
    47 = "ucsb"
}

for ($row = 1; $row -le 47; $row++) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}

# Widen column C to fit the longer school names, as in the edited workbook.
$ws.Columns("C").ColumnWidth = 22.69

# Leave the selection on the last cell that was filled in, matching the
# edited workbook's view state.
$ws.Range("D47").Select()
